$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits in the Skills
#    paragraph (right after "Visual basics").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Rework the "software configuration management tools" sentence: move
#    "Clearcase" so it follows "SVN," instead of following "Serena
#    dimensions,".
# ---------------------------------------------------------------------------
$before = "SVN, Shell scripting, Hudson, Serena dimensions, Clearcase to support code deployment in Apache"
$after  = "SVN, Clearcase and  Shell scripting, Hudson, Serena dimensions to support code deployment in Apache"
$d.Content.Find.Execute($before, $true, $false, $false, $false, $false, $true, 1, $false, $after, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Re-add the "_GoBack" bookmark, now collapsed right after "web logic
#    servers." at the end of that same paragraph (its new location).
# ---------------------------------------------------------------------------
$tail = $d.Content
$tail.Find.Execute("web logic servers.") | Out-Null
$marker = $d.Range($tail.End, $tail.End)
$marker.InsertAfter("~")
$markerRange = $d.Range($tail.End, $tail.End + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerRange2 = $d.Range($tail.End, $tail.End + 1)
$markerRange2.Delete()
